$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 31.14997866666667
$ws.Range("H2").Value = 93.44993600000001
$ws.Range("I2").Value = 0.4621739036316256
$ws.Range("J2").Value = 0.4621739036316256
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07075833333333333
$ws.Range("N2").Value = 0.212275
$ws.Range("O2").Value = 0.006176824525282918
$ws.Range("P2").Value = 0.006176824525282918
$ws.Range("Q2").Value = 2.204120573822222
$ws.Range("R2").Value = 19.8370851644
$ws.Range("S2").Value = 0.002854767102897569
$ws.Range("T2").Value = 0.002854767102897569
# Row 3
$ws.Range("G3").Value = 31.14997866666667
$ws.Range("H3").Value = 93.44993600000001
$ws.Range("I3").Value = 0.4621739036316256
$ws.Range("J3").Value = 0.4621739036316256
$ws.Range("O3").Value = 0.03499706864537662
$ws.Range("P3").Value = 0.03499706864537662
$ws.Range("Q3").Value = 12.48825488064356
$ws.Range("R3").Value = 112.394293925792
$ws.Range("S3").Value = 0.01617473183149768
$ws.Range("T3").Value = 0.01617473183149768
# Row 4
$ws.Range("G4").Value = 31.14997866666667
$ws.Range("H4").Value = 93.44993600000001
$ws.Range("I4").Value = 0.4621739036316256
$ws.Range("J4").Value = 0.4621739036316256
$ws.Range("M4").Value = 0.1538076666666666
$ws.Range("N4").Value = 0.461423
$ws.Range("O4").Value = 0.01342658769487514
$ws.Range("P4").Value = 0.01342658769487514
$ws.Range("Q4").Value = 4.791105535436444
$ws.Range("R4").Value = 43.119949818928
$ws.Range("S4").Value = 0.006205418447392791
$ws.Range("T4").Value = 0.006205418447392791
# Row 5
$ws.Range("G5").Value = 31.14997866666667
$ws.Range("H5").Value = 93.44993600000001
$ws.Range("I5").Value = 0.4621739036316256
$ws.Range("J5").Value = 0.4621739036316256
$ws.Range("M5").Value = 10.82998133333333
$ws.Range("N5").Value = 32.489944
$ws.Range("O5").Value = 0.9453995191344653
$ws.Range("P5").Value = 0.9453995191344652
$ws.Range("Q5").Value = 337.3536874937316
$ws.Range("R5").Value = 3036.183187443584
$ws.Range("S5").Value = 0.4369389862498375
$ws.Range("T5").Value = 0.4369389862498375
# Row 6
$ws.Range("G6").Value = 18.94069966666667
$ws.Range("H6").Value = 56.822099
$ws.Range("I6").Value = 0.2810241764892454
$ws.Range("J6").Value = 0.2810241764892454
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.07075833333333333
$ws.Range("N6").Value = 0.212275
$ws.Range("O6").Value = 0.006176824525282918
$ws.Range("P6").Value = 0.006176824525282918
$ws.Range("Q6").Value = 1.340212340580555
$ws.Range("R6").Value = 12.061911065225
$ws.Range("S6").Value = 0.001735837025536206
$ws.Range("T6").Value = 0.001735837025536206
# Row 7
$ws.Range("G7").Value = 18.94069966666667
$ws.Range("H7").Value = 56.822099
$ws.Range("I7").Value = 0.2810241764892454
$ws.Range("J7").Value = 0.2810241764892454
$ws.Range("O7").Value = 0.03499706864537662
$ws.Range("P7").Value = 0.03499706864537662
$ws.Range("Q7").Value = 7.593465394830889
$ws.Range("R7").Value = 68.34118855347801
$ws.Range("S7").Value = 0.009835022395604555
$ws.Range("T7").Value = 0.009835022395604555
# Row 8
$ws.Range("G8").Value = 18.94069966666667
$ws.Range("H8").Value = 56.822099
$ws.Range("I8").Value = 0.2810241764892454
$ws.Range("J8").Value = 0.2810241764892454
$ws.Range("M8").Value = 0.1538076666666666
$ws.Range("N8").Value = 0.461423
$ws.Range("O8").Value = 0.01342658769487514
$ws.Range("P8").Value = 0.01342658769487514
$ws.Range("Q8").Value = 2.913224820764111
$ws.Range("R8").Value = 26.219023386877
$ws.Range("S8").Value = 0.003773195750012921
$ws.Range("T8").Value = 0.003773195750012921
# Row 9
$ws.Range("G9").Value = 18.94069966666667
$ws.Range("H9").Value = 56.822099
$ws.Range("I9").Value = 0.2810241764892454
$ws.Range("J9").Value = 0.2810241764892454
$ws.Range("M9").Value = 10.82998133333333
$ws.Range("N9").Value = 32.489944
$ws.Range("O9").Value = 0.9453995191344653
$ws.Range("P9").Value = 0.9453995191344652
$ws.Range("Q9").Value = 205.1274238302729
$ws.Range("R9").Value = 1846.146814472456
$ws.Range("S9").Value = 0.2656801213180917
$ws.Range("T9").Value = 0.2656801213180917
# Row 10
$ws.Range("G10").Value = 14.86848
$ws.Range("H10").Value = 44.60544
$ws.Range("I10").Value = 0.2206044349565553
$ws.Range("J10").Value = 0.2206044349565553
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.07075833333333333
$ws.Range("N10").Value = 0.212275
$ws.Range("O10").Value = 0.006176824525282918
$ws.Range("P10").Value = 0.006176824525282918
$ws.Range("Q10").Value = 1.052068864
$ws.Range("R10").Value = 9.468619776000001
$ws.Range("S10").Value = 0.001362634884225831
$ws.Range("T10").Value = 0.001362634884225831
# Row 11
$ws.Range("G11").Value = 14.86848
$ws.Range("H11").Value = 44.60544
$ws.Range("I11").Value = 0.2206044349565553
$ws.Range("J11").Value = 0.2206044349565553
$ws.Range("O11").Value = 0.03499706864537662
$ws.Range("P11").Value = 0.03499706864537662
$ws.Range("Q11").Value = 5.96088266752
$ws.Range("R11").Value = 53.64794400768
$ws.Range("S11").Value = 0.007720508553649088
$ws.Range("T11").Value = 0.007720508553649088
# Row 12
$ws.Range("G12").Value = 14.86848
$ws.Range("H12").Value = 44.60544
$ws.Range("I12").Value = 0.2206044349565553
$ws.Range("J12").Value = 0.2206044349565553
$ws.Range("M12").Value = 0.1538076666666666
$ws.Range("N12").Value = 0.461423
$ws.Range("O12").Value = 0.01342658769487514
$ws.Range("P12").Value = 0.01342658769487514
$ws.Range("Q12").Value = 2.28688621568
$ws.Range("R12").Value = 20.58197594112
$ws.Range("S12").Value = 0.002961964791822568
$ws.Range("T12").Value = 0.002961964791822568
# Row 13
$ws.Range("G13").Value = 14.86848
$ws.Range("H13").Value = 44.60544
$ws.Range("I13").Value = 0.2206044349565553
$ws.Range("J13").Value = 0.2206044349565553
$ws.Range("M13").Value = 10.82998133333333
$ws.Range("N13").Value = 32.489944
$ws.Range("O13").Value = 0.9453995191344653
$ws.Range("P13").Value = 0.9453995191344652
$ws.Range("Q13").Value = 161.02536085504
$ws.Range("R13").Value = 1449.22824769536
$ws.Range("S13").Value = 0.2085593267268578
$ws.Range("T13").Value = 0.2085593267268578
# Row 14
$ws.Range("G14").Value = 2.439668
$ws.Range("H14").Value = 7.319004
$ws.Range("I14").Value = 0.03619748492257375
$ws.Range("J14").Value = 0.03619748492257375
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.07075833333333333
$ws.Range("N14").Value = 0.212275
$ws.Range("O14").Value = 0.006176824525282918
$ws.Range("P14").Value = 0.006176824525282918
$ws.Range("Q14").Value = 0.1726268415666666
$ws.Range("R14").Value = 1.5536415741
$ws.Range("S14").Value = 0.0002235855126233122
$ws.Range("T14").Value = 0.0002235855126233122
# Row 15
$ws.Range("G15").Value = 2.439668
$ws.Range("H15").Value = 7.319004
$ws.Range("I15").Value = 0.03619748492257375
$ws.Range("J15").Value = 0.03619748492257375
$ws.Range("O15").Value = 0.03499706864537662
$ws.Range("P15").Value = 0.03499706864537662
$ws.Range("Q15").Value = 0.9780807920986666
$ws.Range("R15").Value = 8.802727128888
$ws.Range("S15").Value = 0.001266805864625299
$ws.Range("T15").Value = 0.001266805864625299
# Row 16
$ws.Range("G16").Value = 2.439668
$ws.Range("H16").Value = 7.319004
$ws.Range("I16").Value = 0.03619748492257375
$ws.Range("J16").Value = 0.03619748492257375
$ws.Range("M16").Value = 0.1538076666666666
$ws.Range("N16").Value = 0.461423
$ws.Range("O16").Value = 0.01342658769487514
$ws.Range("P16").Value = 0.01342658769487514
$ws.Range("Q16").Value = 0.3752396425213332
$ws.Range("R16").Value = 3.377156782691999
$ws.Range("S16").Value = 0.000486008705646857
$ws.Range("T16").Value = 0.000486008705646857
# Row 17
$ws.Range("G17").Value = 2.439668
$ws.Range("H17").Value = 7.319004
$ws.Range("I17").Value = 0.03619748492257375
$ws.Range("J17").Value = 0.03619748492257375
$ws.Range("M17").Value = 10.82998133333333
$ws.Range("N17").Value = 32.489944
$ws.Range("O17").Value = 0.9453995191344653
$ws.Range("P17").Value = 0.9453995191344652
$ws.Range("Q17").Value = 26.42155889953067
$ws.Range("R17").Value = 237.794030095776
$ws.Range("S17").Value = 0.03422108483967828
$ws.Range("T17").Value = 0.03422108483967828
